# Adds two new paragraphs of body text around the existing blank paragraph
# that sits right before the final section break, per the target diff:
#   ... <blank para> / <NEW para "One of the biggest issues..."> / <blank para/> /
#   <NEW para "Programmers are often..."> / <sectPr>
#
# Note: paragraph objects captured before a mutation go stale afterwards
# (their .Previous/.Next no longer resolve), so every paragraph we touch is
# re-fetched fresh by index from $d.Paragraphs right before use.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Locate the lone blank paragraph that immediately precedes the
#     section break (the very last paragraph in the main story). ---
$anchorIndex = $d.Paragraphs.Count
$anchorStart = $d.Paragraphs.Item($anchorIndex).Range.Start

# --- Paragraph 1: inserted BEFORE the anchor paragraph ---
$insBefore = $d.Range($anchorStart, $anchorStart)
$null = $insBefore.InsertParagraphBefore()

# The freshly inserted empty paragraph now occupies the anchor's old slot;
# the (untouched) anchor paragraph shifted one slot later.
$firstRange = $d.Paragraphs.Item($anchorIndex).Range
$firstXml = '<w:p ' + $wNs + '>' `
  + '<w:r><w:t xml:space="preserve">One of the biggest issues that plagues the programming profession is the fragility of projects created and worked on local machines. Projects interfacing with software components of computers are given unparallel </w:t></w:r>' `
  + '<w:r><w:t xml:space="preserve">access to its functionality, with that, comes a lack of security and stability. The ability to create and design with computer capabilities has been and always will be utilised endlessly, causing the need for a solution that includes preventative or </w:t></w:r>' `
  + '<w:r><w:t>corrective measures that at the very least allow for the salvaging of any completed work.</w:t></w:r>' `
  + '</w:p>'
$null = $firstRange.InsertXML($firstXml)

$anchorIndex = $anchorIndex + 1

# --- Paragraph 2: inserted AFTER the anchor paragraph ---
$anchorEnd = $d.Paragraphs.Item($anchorIndex).Range.End
$insAfter = $d.Range($anchorEnd, $anchorEnd)
$null = $insAfter.InsertParagraphAfter()

$secondRange = $d.Paragraphs.Item($anchorIndex + 1).Range
$secondXml = '<w:p ' + $wNs + '>' `
  + '<w:r><w:t xml:space="preserve">Programmers are often </w:t></w:r>' `
  + '<w:r><w:t xml:space="preserve">the ones at the end of the stick, Creating large changes, providing new computer instructions, accessing and combing through the inner workings of the computer which has the high potential to bump something accidentally causing unforeseen problems to occur. </w:t></w:r>' `
  + '</w:p>'
$null = $secondRange.InsertXML($secondXml)

Write-Output "Paragraphs after edit: $($d.Paragraphs.Count)"
